# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Writes the newly calculated s_vals (K) for rows 2-35 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 3
    4  = 6
    5  = 7
    6  = 8
    7  = 7
    8  = 6
    9  = 7
    10 = 8
    11 = 3
    12 = 6
    13 = 2
    14 = 4
    15 = 5
    16 = 5
    17 = 7
    18 = 4
    19 = 8
    20 = 3
    21 = 6
    22 = 7
    23 = 8
    24 = 4
    25 = 7
    26 = 5
    27 = 5
    28 = 4
    29 = 8
    30 = 3
    31 = 5
    32 = 4
    33 = 3
    34 = 7
    35 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
